$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Update marking points per correct answer (3 -> 5)
$ws.Range("B11").Value = 5

# Update total obtained marks (66 -> 110) = 22 correct * 5
$ws.Range("B12").Value = 110

# Update the "obtained/max" label (66/84 -> 110/140) = 110 obtained / (28 max questions * 5)
$ws.Range("E12").Value = "110/140"
